$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.184.19"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.022.21"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "3.013.50"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000217"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "3.521.82"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "62.360.78"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "3.030.96"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "58.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "464.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "
$ws.Range("D38").Value = "3.186.72"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0384"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0776"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.244"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.81%  "
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.81%  "
